# Auto-generated edit script: updates TPM-derived LR-pair metrics (Col4a1-Itgav)
# per commit "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 338.6742503333333
$ws.Cells.Item(2, 8).Value = 1016.022751
$ws.Cells.Item(2, 9).Value = 0.5849329800180821
$ws.Cells.Item(2, 10).Value = 0.584932980018082
$ws.Cells.Item(2, 13).Value = 8.820647333333334
$ws.Cells.Item(2, 14).Value = 26.461942
$ws.Cells.Item(2, 15).Value = 0.06415146660411865
$ws.Cells.Item(2, 16).Value = 0.06415146660411865
$ws.Cells.Item(2, 17).Value = 2987.326123071382
$ws.Cells.Item(2, 18).Value = 26885.93510764244
$ws.Cells.Item(2, 19).Value = 0.03752430853327759
$ws.Cells.Item(2, 20).Value = 0.03752430853327759
$ws.Cells.Item(3, 7).Value = 338.6742503333333
$ws.Cells.Item(3, 8).Value = 1016.022751
$ws.Cells.Item(3, 9).Value = 0.5849329800180821
$ws.Cells.Item(3, 10).Value = 0.584932980018082
$ws.Cells.Item(3, 15).Value = 0.3979101621202897
$ws.Cells.Item(3, 16).Value = 0.3979101621202898
$ws.Cells.Item(3, 17).Value = 18529.38810071093
$ws.Cells.Item(3, 18).Value = 166764.4929063983
$ws.Cells.Item(3, 19).Value = 0.2327507769084993
$ws.Cells.Item(3, 20).Value = 0.2327507769084992
$ws.Cells.Item(4, 7).Value = 338.6742503333333
$ws.Cells.Item(4, 8).Value = 1016.022751
$ws.Cells.Item(4, 9).Value = 0.5849329800180821
$ws.Cells.Item(4, 10).Value = 0.584932980018082
$ws.Cells.Item(4, 13).Value = 21.90816333333333
$ws.Cells.Item(4, 14).Value = 65.72449
$ws.Cells.Item(4, 15).Value = 0.1593353362087987
$ws.Cells.Item(4, 16).Value = 0.1593353362087987
$ws.Cells.Item(4, 17).Value = 7419.730793096887
$ws.Cells.Item(4, 18).Value = 66777.57713787199
$ws.Cells.Item(4, 19).Value = 0.09320049303079563
$ws.Cells.Item(4, 20).Value = 0.09320049303079561
$ws.Cells.Item(5, 7).Value = 338.6742503333333
$ws.Cells.Item(5, 8).Value = 1016.022751
$ws.Cells.Item(5, 9).Value = 0.5849329800180821
$ws.Cells.Item(5, 10).Value = 0.584932980018082
$ws.Cells.Item(5, 13).Value = 52.056859
$ws.Cells.Item(5, 14).Value = 156.170577
$ws.Cells.Item(5, 15).Value = 0.3786030350667928
$ws.Cells.Item(5, 16).Value = 0.3786030350667929
$ws.Cells.Item(5, 17).Value = 17630.31769653303
$ws.Cells.Item(5, 18).Value = 158672.8592687973
$ws.Cells.Item(5, 19).Value = 0.2214574015455096
$ws.Cells.Item(5, 20).Value = 0.2214574015455096
$ws.Cells.Item(6, 9).Value = 0.279688040971731
$ws.Cells.Item(6, 10).Value = 0.2796880409717309
$ws.Cells.Item(6, 13).Value = 8.820647333333334
$ws.Cells.Item(6, 14).Value = 26.461942
$ws.Cells.Item(6, 15).Value = 0.06415146660411865
$ws.Cells.Item(6, 16).Value = 0.06415146660411865
$ws.Cells.Item(6, 17).Value = 1428.401918933828
$ws.Cells.Item(6, 18).Value = 12855.61727040445
$ws.Cells.Item(6, 19).Value = 0.01794239801996937
$ws.Cells.Item(6, 20).Value = 0.01794239801996936
$ws.Cells.Item(7, 9).Value = 0.279688040971731
$ws.Cells.Item(7, 10).Value = 0.2796880409717309
$ws.Cells.Item(7, 15).Value = 0.3979101621202897
$ws.Cells.Item(7, 16).Value = 0.3979101621202898
$ws.Cells.Item(7, 19).Value = 0.1112907137261677
$ws.Cells.Item(7, 20).Value = 0.1112907137261677
$ws.Cells.Item(8, 9).Value = 0.279688040971731
$ws.Cells.Item(8, 10).Value = 0.2796880409717309
$ws.Cells.Item(8, 13).Value = 21.90816333333333
$ws.Cells.Item(8, 14).Value = 65.72449
$ws.Cells.Item(8, 15).Value = 0.1593353362087987
$ws.Cells.Item(8, 16).Value = 0.1593353362087987
$ws.Cells.Item(8, 17).Value = 3547.773917611459
$ws.Cells.Item(8, 18).Value = 31929.96525850313
$ws.Cells.Item(8, 19).Value = 0.04456418804181102
$ws.Cells.Item(8, 20).Value = 0.044564188041811
$ws.Cells.Item(9, 9).Value = 0.279688040971731
$ws.Cells.Item(9, 10).Value = 0.2796880409717309
$ws.Cells.Item(9, 13).Value = 52.056859
$ws.Cells.Item(9, 14).Value = 156.170577
$ws.Cells.Item(9, 15).Value = 0.3786030350667928
$ws.Cells.Item(9, 16).Value = 0.3786030350667929
$ws.Cells.Item(9, 17).Value = 8430.006832748826
$ws.Cells.Item(9, 18).Value = 75870.06149473944
$ws.Cells.Item(9, 19).Value = 0.1058907411837829
$ws.Cells.Item(9, 20).Value = 0.1058907411837829
$ws.Cells.Item(10, 7).Value = 77.79536166666666
$ws.Cells.Item(10, 8).Value = 233.386085
$ws.Cells.Item(10, 9).Value = 0.1343623634996766
$ws.Cells.Item(10, 10).Value = 0.1343623634996766
$ws.Cells.Item(10, 13).Value = 8.820647333333334
$ws.Cells.Item(10, 14).Value = 26.461942
$ws.Cells.Item(10, 15).Value = 0.06415146660411865
$ws.Cells.Item(10, 16).Value = 0.06415146660411865
$ws.Cells.Item(10, 17).Value = 686.2054494307855
$ws.Cells.Item(10, 18).Value = 6175.84904487707
$ws.Cells.Item(10, 19).Value = 0.008619542674899955
$ws.Cells.Item(10, 20).Value = 0.008619542674899953
$ws.Cells.Item(11, 7).Value = 77.79536166666666
$ws.Cells.Item(11, 8).Value = 233.386085
$ws.Cells.Item(11, 9).Value = 0.1343623634996766
$ws.Cells.Item(11, 10).Value = 0.1343623634996766
$ws.Cells.Item(11, 15).Value = 0.3979101621202897
$ws.Cells.Item(11, 16).Value = 0.3979101621202898
$ws.Cells.Item(11, 17).Value = 4256.303652663491
$ws.Cells.Item(11, 18).Value = 38306.73287397143
$ws.Cells.Item(11, 19).Value = 0.05346414984302162
$ws.Cells.Item(11, 20).Value = 0.05346414984302162
$ws.Cells.Item(12, 7).Value = 77.79536166666666
$ws.Cells.Item(12, 8).Value = 233.386085
$ws.Cells.Item(12, 9).Value = 0.1343623634996766
$ws.Cells.Item(12, 10).Value = 0.1343623634996766
$ws.Cells.Item(12, 13).Value = 21.90816333333333
$ws.Cells.Item(12, 14).Value = 65.72449
$ws.Cells.Item(12, 15).Value = 0.1593353362087987
$ws.Cells.Item(12, 16).Value = 0.1593353362087987
$ws.Cells.Item(12, 17).Value = 1704.353489969072
$ws.Cells.Item(12, 18).Value = 15339.18140972165
$ws.Cells.Item(12, 19).Value = 0.02140867236202979
$ws.Cells.Item(12, 20).Value = 0.02140867236202979
$ws.Cells.Item(13, 7).Value = 77.79536166666666
$ws.Cells.Item(13, 8).Value = 233.386085
$ws.Cells.Item(13, 9).Value = 0.1343623634996766
$ws.Cells.Item(13, 10).Value = 0.1343623634996766
$ws.Cells.Item(13, 13).Value = 52.056859
$ws.Cells.Item(13, 14).Value = 156.170577
$ws.Cells.Item(13, 15).Value = 0.3786030350667928
$ws.Cells.Item(13, 16).Value = 0.3786030350667929
$ws.Cells.Item(13, 17).Value = 4049.782173135671
$ws.Cells.Item(13, 18).Value = 36448.03955822104
$ws.Cells.Item(13, 19).Value = 0.05086999861972523
$ws.Cells.Item(13, 20).Value = 0.05086999861972523
$ws.Cells.Item(14, 7).Value = 0.5886170000000001
$ws.Cells.Item(14, 8).Value = 1.765851
$ws.Cells.Item(14, 9).Value = 0.001016615510510267
$ws.Cells.Item(14, 10).Value = 0.001016615510510266
$ws.Cells.Item(14, 13).Value = 8.820647333333334
$ws.Cells.Item(14, 14).Value = 26.461942
$ws.Cells.Item(14, 15).Value = 0.06415146660411865
$ws.Cells.Item(14, 16).Value = 0.06415146660411865
$ws.Cells.Item(14, 17).Value = 5.191982971404667
$ws.Cells.Item(14, 18).Value = 46.727846742642
$ws.Cells.Item(14, 19).Value = 0.0000652173759717284
$ws.Cells.Item(14, 20).Value = 0.00006521737597172838
$ws.Cells.Item(15, 7).Value = 0.5886170000000001
$ws.Cells.Item(15, 8).Value = 1.765851
$ws.Cells.Item(15, 9).Value = 0.001016615510510267
$ws.Cells.Item(15, 10).Value = 0.001016615510510266
$ws.Cells.Item(15, 15).Value = 0.3979101621202897
$ws.Cells.Item(15, 16).Value = 0.3979101621202898
$ws.Cells.Item(15, 17).Value = 32.204139597095
$ws.Cells.Item(15, 18).Value = 289.837256373855
$ws.Cells.Item(15, 19).Value = 0.0004045216426011413
$ws.Cells.Item(15, 20).Value = 0.0004045216426011413
$ws.Cells.Item(16, 7).Value = 0.5886170000000001
$ws.Cells.Item(16, 8).Value = 1.765851
$ws.Cells.Item(16, 9).Value = 0.001016615510510267
$ws.Cells.Item(16, 10).Value = 0.001016615510510266
$ws.Cells.Item(16, 13).Value = 21.90816333333333
$ws.Cells.Item(16, 14).Value = 65.72449
$ws.Cells.Item(16, 15).Value = 0.1593353362087987
$ws.Cells.Item(16, 16).Value = 0.1593353362087987
$ws.Cells.Item(16, 17).Value = 12.89551737677667
$ws.Cells.Item(16, 18).Value = 116.05965639099
$ws.Cells.Item(16, 19).Value = 0.0001619827741622328
$ws.Cells.Item(16, 20).Value = 0.0001619827741622328
$ws.Cells.Item(17, 7).Value = 0.5886170000000001
$ws.Cells.Item(17, 8).Value = 1.765851
$ws.Cells.Item(17, 9).Value = 0.001016615510510267
$ws.Cells.Item(17, 10).Value = 0.001016615510510266
$ws.Cells.Item(17, 13).Value = 52.056859
$ws.Cells.Item(17, 14).Value = 156.170577
$ws.Cells.Item(17, 15).Value = 0.3786030350667928
$ws.Cells.Item(17, 16).Value = 0.3786030350667929
$ws.Cells.Item(17, 17).Value = 30.641552174003
$ws.Cells.Item(17, 18).Value = 275.7739695660269
$ws.Cells.Item(17, 19).Value = 0.000384893717775164
$ws.Cells.Item(17, 20).Value = 0.000384893717775164
